# Late night intro edits
$d = $word.ActiveDocument

# 1. Second figure caption: "Figure 1" -> "Figure 2"
#    (the paragraph that begins "Maternal telomere length and offspring
#    condition and survival..."). Scope the Find to that paragraph only so
#    the first figure's "Figure 1" label is left untouched.
$p3 = $d.Paragraphs(3).Range
$p3.Find.ClearFormatting()
$p3.Find.Execute("Figure 1", $true, $false, $false, $false, $false, $true, 1, $false, "Figure 2", 2)

# 2. Panel B caption sentence: drop "maternal" and switch "logistic" -> "linear"
$d.Content.Find.ClearFormatting()
$d.Content.Find.Execute(
    "Offspring survival to adulthood in relation to maternal telomere length. Lines and shaded areas represent estimates and 95% confidence limits from a logistic regression (",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Offspring survival to adulthood in relation to telomere length. Lines and shaded areas represent estimates and 95% confidence limits from a linear regression (",
    2)

# 3. Swap the regression type named after panel (A)
$d.Content.Find.ClearFormatting()
$d.Content.Find.Execute(
    ") linear regression (",
    $true, $false, $false, $false, $false, $true, 1, $false,
    ") and logistic regression (",
    2)
